$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.579.45"
$ws.Range("E2").Value = "  +2.19%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.914.52"
$ws.Range("E3").Value = "  +5.64%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.12%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'315.69"
$ws.Range("E5").Value = "  +1.88%  "

# Row 6 - USDC
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.15%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.5165"
$ws.Range("E7").Value = "  +3.55%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3981"
$ws.Range("E8").Value = "  +1.62%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.09761"
$ws.Range("E9").Value = "  -1.81%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "'1.159"
$ws.Range("E10").Value = "  +5.31%  "

# Row 11 - OKB
$ws.Range("D11").Value = "'42.08"
$ws.Range("E11").Value = "  +2.91%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "'6.552"
$ws.Range("E12").Value = "  +2.27%  "

# Row 13 - Solana
$ws.Range("D13").Value = "'21.29"
$ws.Range("E13").Value = "  +3.94%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.914.42"
$ws.Range("E14").Value = "  +5.85%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "'7.597"
$ws.Range("E15").Value = "  +4.68%  "

# Row 16 - BinanceUSD
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.14%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "'0.00001142"
$ws.Range("E17").Value = "  -0.30%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "'94.06"

# Row 19 - TRON
$ws.Range("D19").Value = "'0.06656"
$ws.Range("E19").Value = "  +0.13%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "'18.17"
$ws.Range("E20").Value = "  +5.87%  "

# Row 21 - Dai
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.01%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +6.93%  "

# Row 23 - WrappedBTC
$ws.Range("D23").Value = "28.645.80"
$ws.Range("E23").Value = "  +2.21%  "

# Row 24 - Cosmos
$ws.Range("D24").Value = "'11.49"
$ws.Range("E24").Value = "  +3.87%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +1.66%  "

# Row 26 - LEO (was LidoDAOToken)
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "'3.394"
$ws.Range("E26").Value = "  -0.72%  "

# Row 27 - LidoDAOToken (was LEO)
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.698"
$ws.Range("E27").Value = "  +12.62%  "

# Row 28 - WrappedliquidstakedEther2.0
$ws.Range("D28").Value = "2.133.29"
$ws.Range("E28").Value = "  +5.64%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "'21.28"
$ws.Range("E29").Value = "  +3.53%  "

# Row 30 - Monero
$ws.Range("D30").Value = "'160.22"
$ws.Range("E30").Value = "  +1.16%  "

# Row 31 - BitcoinCash
$ws.Range("D31").Value = "'129.09"
$ws.Range("E31").Value = "  +1.70%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "'1.108"
$ws.Range("E32").Value = "  +7.60%  "

# Row 33 - Stellar
$ws.Range("D33").Value = "'0.1084"
$ws.Range("E33").Value = "  +2.15%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "'5.770"
$ws.Range("E34").Value = "  +3.58%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "'3.641"
$ws.Range("E35").Value = "  +0.94%  "

# Row 36 - FraxShare
$ws.Range("D36").Value = "'9.878"
$ws.Range("E36").Value = "  +11.35%  "

# Row 37 - Hedera
$ws.Range("D37").Value = "'0.06809"
$ws.Range("E37").Value = "  +1.34%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "'0.02440"
$ws.Range("E38").Value = "  +4.98%  "

# Row 39 - ARBITRUM
$ws.Range("D39").Value = "'1.270"
$ws.Range("E39").Value = "  +8.26%  "

# Row 40 - Algorand
$ws.Range("D40").Value = "'0.2233"
$ws.Range("E40").Value = "  +4.36%  "

# Row 41 - Aptos
$ws.Range("D41").Value = "'11.89"
$ws.Range("E41").Value = "  +5.47%  "

# Row 42 - InternetComputer(DFINITY)
$ws.Range("D42").Value = "'5.134"
$ws.Range("E42").Value = "  +4.05%  "

# Row 43 - TheSandbox
$ws.Range("D43").Value = "'0.6446"
$ws.Range("E43").Value = "  +4.37%  "

# Row 44 - TrustWalletToken
$ws.Range("E44").Value = "  +1.91%  "

# Row 45 - Frax
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  +0.07%  "

# Row 46 - EnergySwap
$ws.Range("D46").Value = "'13.64"
$ws.Range("E46").Value = "  +4.02%  "

# Row 47 - Decentraland
$ws.Range("D47").Value = "'0.6091"
$ws.Range("E47").Value = "  +3.41%  "

# Row 48 - PancakeSwap
$ws.Range("D48").Value = "'3.789"
$ws.Range("E48").Value = "  +2.76%  "

# Row 49 - WEMIXTOKEN
$ws.Range("D49").Value = "'1.281"
$ws.Range("E49").Value = "  +0.15%  "

# Row 50 - NEARProtocol
$ws.Range("E50").Value = "  +5.83%  "

# Row 51 - Quant
$ws.Range("D51").Value = "'125.33"
$ws.Range("E51").Value = "  +1.36%  "
